$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.6091841905957551
$ws.Range("C2").Value = 0.07670933973429328
$ws.Range("D2").Value = 0.07823536744989212
$ws.Range("E2").Value = 0.4092159168139062
$ws.Range("G2").Value = 1.189793201367877
$ws.Range("H2").Value = 1.101868511943223
$ws.Range("K2").Value = 0.6090929765921942
$ws.Range("N2").Value = 1.867259638893344
$ws.Range("B3").Value = 0.5562926210567696
$ws.Range("C3").Value = 0.06680845881504638
$ws.Range("D3").Value = 0.07099624596710896
$ws.Range("E3").Value = 0.356956660429546
$ws.Range("G3").Value = 1.164746019747582
$ws.Range("H3").Value = 1.095636539561326
$ws.Range("K3").Value = 0.548940604643974
$ws.Range("N3").Value = 1.877060035399715
$ws.Range("B4").Value = 0.5241275071028895
$ws.Range("C4").Value = 0.06073559980339382
$ws.Range("D4").Value = 0.06658999394001341
$ws.Range("E4").Value = 0.3249800277407502
$ws.Range("G4").Value = 1.150107554988253
$ws.Range("H4").Value = 1.092353102637247
$ws.Range("K4").Value = 0.5122801669463399
$ws.Range("N4").Value = 1.88368721243954
$ws.Range("B5").Value = 0.5110974694362937
$ws.Range("C5").Value = 0.05826225952571917
$ws.Range("D5").Value = 0.06480401988150675
$ws.Range("E5").Value = 0.3119749038452255
$ws.Range("G5").Value = 1.144327410036041
$ws.Range("H5").Value = 1.091151252381536
$ws.Range("K5").Value = 0.4974087054447125
$ws.Range("N5").Value = 1.886540701828004
$ws.Range("B6").Value = 0.5089385062962037
$ws.Range("C6").Value = 0.05785164223715356
$ws.Range("D6").Value = 0.0645080378497056
$ws.Range("E6").Value = 0.3098168948926201
$ws.Range("G6").Value = 1.143378772414991
$ws.Range("H6").Value = 1.090959900121987
$ws.Range("K6").Value = 0.4949433931485032
$ws.Range("N6").Value = 1.887023740866695
$ws.Range("B7").Value = 0.5239514663880698
$ws.Range("C7").Value = 0.06070223805618014
$ws.Range("D7").Value = 0.06656586888668414
$ws.Range("E7").Value = 0.3248045353079618
$ws.Range("G7").Value = 1.150028853712413
$ws.Range("H7").Value = 1.092336343189729
$ws.Range("K7").Value = 0.5120793310436511
$ws.Range("N7").Value = 1.883725077284083
$ws.Range("B8").Value = 0.59088236129287
$ws.Range("C8").Value = 0.07329408569594875
$ws.Range("D8").Value = 0.07573124682662069
$ws.Range("E8").Value = 0.3911725693668586
$ws.Range("G8").Value = 1.181002555315473
$ws.Range("H8").Value = 1.099606795678227
$ws.Range("K8").Value = 0.5882952139473048
$ws.Range("N8").Value = 1.87051200049271
$ws.Range("B9").Value = 0.7246310723945726
$ws.Range("C9").Value = 0.09804666490794034
$ws.Range("D9").Value = 0.09401611826045553
$ws.Range("E9").Value = 0.5223063258812601
$ws.Range("G9").Value = 1.247671898878707
$ws.Range("H9").Value = 1.118191608240409
$ws.Range("K9").Value = 0.7399666090250037
$ws.Range("N9").Value = 1.849459487183026
$ws.Range("B10").Value = 0.8244761546641826
$ws.Range("C10").Value = 0.1162843902163786
$ws.Range("D10").Value = 0.1076487884881914
$ws.Range("E10").Value = 0.6194095504877453
$ws.Range("G10").Value = 1.300350118332574
$ws.Range("H10").Value = 1.134513302833568
$ws.Range("K10").Value = 0.8528202076288096
$ws.Range("N10").Value = 1.836982933797742
$ws.Range("B11").Value = 0.870254539771679
$ws.Range("C11").Value = 0.1245959662509506
$ws.Range("D11").Value = 0.1138958088410931
$ws.Range("E11").Value = 0.6637841985572805
$ws.Range("G11").Value = 1.32513551884702
$ws.Range("H11").Value = 1.142524297557259
$ws.Range("K11").Value = 0.9044851475094902
$ws.Range("N11").Value = 1.831962831121885
$ws.Range("B12").Value = 0.8876420153180788
$ws.Range("C12").Value = 0.1277457691148811
$ws.Range("D12").Value = 0.1162680621956866
$ws.Range("E12").Value = 0.6806194584708152
$ws.Range("G12").Value = 1.334640638273441
$ws.Range("H12").Value = 1.145642638202929
$ws.Range("K12").Value = 0.9240974794723513
$ws.Range("N12").Value = 1.830156654648192
$ws.Range("B13").Value = 0.8838949792648236
$ws.Range("C13").Value = 0.1270672930332921
$ws.Range("D13").Value = 0.1157568578573347
$ws.Range("E13").Value = 0.6769922407611944
$ws.Range("G13").Value = 1.332588209139203
$ws.Range("H13").Value = 1.144967271020846
$ws.Range("K13").Value = 0.9198714665202772
$ws.Range("N13").Value = 1.830541420565069
$ws.Range("B14").Value = 0.8716839679385089
$ws.Range("C14").Value = 0.1248550533152013
$ws.Range("D14").Value = 0.114090841814118
$ws.Range("E14").Value = 0.6651685988295952
$ws.Range("G14").Value = 1.325915110045685
$ws.Range("H14").Value = 1.142779144124546
$ws.Range("K14").Value = 0.906097699965386
$ws.Range("N14").Value = 1.83181233217158
$ws.Range("B15").Value = 0.8642111881386541
$ws.Range("C15").Value = 0.1235003101279801
$ws.Range("D15").Value = 0.1130712270997662
$ws.Range("E15").Value = 0.6579304597945139
$ws.Range("G15").Value = 1.321843234364621
$ws.Range("H15").Value = 1.141449904060721
$ws.Range("K15").Value = 0.8976671399043425
$ws.Range("N15").Value = 1.832603168348385
$ws.Range("B16").Value = 0.8214916969589581
$ws.Range("C16").Value = 0.1157415319263464
$ws.Range("D16").Value = 0.1072414546719784
$ws.Range("E16").Value = 0.6165138394653553
$ws.Range("G16").Value = 1.298746977629605
$ws.Range("H16").Value = 1.1340015999115
$ws.Range("K16").Value = 0.8494504426099354
$ws.Range("N16").Value = 1.837324248306857
$ws.Range("B17").Value = 0.7953769238727091
$ws.Range("C17").Value = 0.1109858010116227
$ws.Range("D17").Value = 0.1036767944062404
$ws.Range("E17").Value = 0.5911595716423363
$ws.Range("G17").Value = 1.284789543026733
$ws.Range("H17").Value = 1.12958276420278
$ws.Range("K17").Value = 0.8199555028149632
$ws.Range("N17").Value = 1.840388775086495
$ws.Range("B18").Value = 0.7803900867390894
$ws.Range("C18").Value = 0.1082518428602839
$ws.Range("D18").Value = 0.1016307676299704
$ws.Range("E18").Value = 0.5765953192967714
$ws.Range("G18").Value = 1.276838873132562
$ws.Range("H18").Value = 1.127096312264598
$ws.Range("K18").Value = 0.8030215442964277
$ws.Range("N18").Value = 1.842213078724114
$ws.Range("B19").Value = 0.7753215674479179
$ws.Range("C19").Value = 0.1073264078059992
$ws.Range("D19").Value = 0.1009387485348583
$ws.Range("E19").Value = 0.5716672756463907
$ws.Range("G19").Value = 1.274160147532797
$ws.Range("H19").Value = 1.126263899772397
$ws.Range("K19").Value = 0.7972932518736116
$ws.Range("N19").Value = 1.842841331304825
$ws.Range("B20").Value = 0.7981533945392698
$ws.Range("C20").Value = 0.1114919089710042
$ws.Range("D20").Value = 0.1040558157195193
$ws.Range("E20").Value = 0.5938566075170968
$ws.Range("G20").Value = 1.286267328876505
$ws.Range("H20").Value = 1.130047446185102
$ws.Range("K20").Value = 0.8230921006031338
$ws.Range("N20").Value = 1.840056164487848
$ws.Range("B21").Value = 0.875269216056779
$ws.Range("C21").Value = 0.1255047751926668
$ws.Range("D21").Value = 0.1145800099821201
$ws.Range("E21").Value = 0.6686406128029319
$ws.Range("G21").Value = 1.327871909769982
$ws.Range("H21").Value = 1.143419546215739
$ws.Range("K21").Value = 0.9101420808979412
$ws.Range("N21").Value = 1.831436456522439
$ws.Range("B22").Value = 0.9259734120518033
$ws.Range("C22").Value = 0.1346770213670823
$ws.Range("D22").Value = 0.1214969464745224
$ws.Range("E22").Value = 0.7177013354224613
$ws.Range("G22").Value = 1.355759607507821
$ws.Range("H22").Value = 1.152653198556521
$ws.Range("K22").Value = 0.9673141297860184
$ws.Range("N22").Value = 1.826355972285626
$ws.Range("B23").Value = 0.8988835391048156
$ws.Range("C23").Value = 0.1297802665133077
$ws.Range("D23").Value = 0.1178016644558824
$ws.Range("E23").Value = 0.6914989333135964
$ws.Range("G23").Value = 1.340811252982974
$ws.Range("H23").Value = 1.147679649988163
$ws.Range("K23").Value = 0.9367744295189766
$ws.Range("N23").Value = 1.829016735907487
$ws.Range("B24").Value = 0.7968980684334497
$ws.Range("C24").Value = 0.1112630970351063
$ws.Range("D24").Value = 0.1038844497781639
$ws.Range("E24").Value = 0.5926372396447874
$ws.Range("G24").Value = 1.285598992587779
$ws.Range("H24").Value = 1.129837195296233
$ws.Range("K24").Value = 0.821673972938413
$ws.Range("N24").Value = 1.840206343191369
$ws.Range("B25").Value = 0.6881752520068005
$ws.Range("C25").Value = 0.09134242952094951
$ws.Range("D25").Value = 0.08903531152128608
$ws.Range("E25").Value = 0.4867093695847018
$ws.Range("G25").Value = 1.22899289237526
$ws.Range("H25").Value = 1.112697639411437
$ws.Range("K25").Value = 0.6986912594324792
$ws.Range("N25").Value = 1.854631668916639
